# Spirit karma_performance.xlsx rework:
#  - rename the three sheets
#  - split the two benchmark tables that used to live on Sheet1 across
#    "Single double" (unchanged rows 3-10) and "Sequence of items" (was rows
#    36-42, now rows 3-9 of its own sheet)
#  - add a new "Single int" benchmark table + chart on the third sheet
#  - move the "sequence of items" chart off Sheet1 onto its new sheet, and
#    add a new chart for the int benchmark
#  - misc selection / active-tab bookkeeping to match the saved-from-Excel
#    state

$wb = $excel.ActiveWorkbook

$wsDouble = $wb.Worksheets.Item(1)
$wsSeq    = $wb.Worksheets.Item(2)
$wsInt    = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# 1. Move the "sequence of items" table (old D36:P42) from sheet1 to the
#    (still empty) second sheet, landing at D3:P9.
# ---------------------------------------------------------------------
$wsDouble.Range("D36:P42").Copy()
$wsSeq.Range("D3").PasteSpecial()
$excel.CutCopyMode = $false

# match column widths of the source table
$wsSeq.Columns.Item(4).ColumnWidth  = $wsDouble.Columns.Item(4).ColumnWidth
$wsSeq.Columns.Item(5).ColumnWidth  = $wsDouble.Columns.Item(5).ColumnWidth
$wsSeq.Columns.Item(6).ColumnWidth  = $wsDouble.Columns.Item(6).ColumnWidth
$wsSeq.Columns.Item(7).ColumnWidth  = $wsDouble.Columns.Item(7).ColumnWidth
$wsSeq.Columns.Item(8).ColumnWidth  = $wsDouble.Columns.Item(8).ColumnWidth
$wsSeq.Columns.Item(9).ColumnWidth  = $wsDouble.Columns.Item(9).ColumnWidth
$wsSeq.Columns.Item(10).ColumnWidth = $wsDouble.Columns.Item(10).ColumnWidth
$wsSeq.Columns.Item(15).ColumnWidth = $wsDouble.Columns.Item(15).ColumnWidth
$wsSeq.Columns.Item(16).ColumnWidth = $wsDouble.Columns.Item(16).ColumnWidth

# now remove the old copy (rows 36-42) from sheet1 and shrink its col D
$wsDouble.Range("D36:P42").ClearContents()
$wsDouble.Columns.Item(4).ColumnWidth = 12.88

# ---------------------------------------------------------------------
# 2. Populate the new "Single int" benchmark table on sheet3.
# ---------------------------------------------------------------------
$wsInt.Range("D1").Value = "Converting 10000000 randomly generated int values to strings."

$wsInt.Range("E3").Value = $wsDouble.Range("E3").Value
$wsInt.Range("F3").Value = $wsDouble.Range("F3").Value
$wsInt.Range("G3").Value = $wsDouble.Range("G3").Value
$wsInt.Range("H3").Value = $wsDouble.Range("H3").Value
$wsInt.Range("I3").Value = $wsDouble.Range("I3").Value
$wsInt.Range("J3").Value = $wsDouble.Range("J3").Value

$wsInt.Range("D4").Value = "ltoa "
$wsInt.Range("E4").Value = 1.5589999999999999
$wsInt.Range("F4").Value = 0.88100000000000001
$wsInt.Range("G4").Value = 0.85899999999999999
$wsInt.Range("H4").Value = 1.2
$wsInt.Range("I4").Value = 1.1180000000000001
$wsInt.Range("J4").Value = 0.88900000000000001

$wsInt.Range("D5").Value = "iostreams   "
$wsInt.Range("E5").Value = 6.484
$wsInt.Range("F5").Value = 13.161
$wsInt.Range("G5").Value = 11.635999999999999
$wsInt.Range("H5").Value = 3.42
$wsInt.Range("I5").Value = 7.8120000000000003
$wsInt.Range("J5").Value = 7.3680000000000003

$wsInt.Range("D6").Value = "Boost.Format"
$wsInt.Range("E6").Value = 16.823
$wsInt.Range("F6").Value = 21.568999999999999
$wsInt.Range("G6").Value = 19.706
$wsInt.Range("H6").Value = 17.28
$wsInt.Range("I6").Value = 14.401999999999999
$wsInt.Range("J6").Value = 13.222

$wsInt.Range("D7").Value = "Karma"
$wsInt.Range("E7").Value = 2.5619999999999998
$wsInt.Range("F7").Value = 1.0109999999999999
$wsInt.Range("G7").Value = 0.95499999999999996
$wsInt.Range("H7").Value = 2.956
$wsInt.Range("I7").Value = 1.016
$wsInt.Range("J7").Value = 0.878

# E12 carries a wrap-text style with no content (left behind by the author)
$wsInt.Range("E12").WrapText = $true

$wsInt.Columns.Item(4).ColumnWidth  = $wsSeq.Columns.Item(4).ColumnWidth
$wsInt.Columns.Item(5).ColumnWidth  = $wsSeq.Columns.Item(5).ColumnWidth
$wsInt.Columns.Item(6).ColumnWidth  = $wsSeq.Columns.Item(6).ColumnWidth
$wsInt.Columns.Item(7).ColumnWidth  = $wsSeq.Columns.Item(7).ColumnWidth
$wsInt.Columns.Item(8).ColumnWidth  = $wsSeq.Columns.Item(8).ColumnWidth
$wsInt.Columns.Item(9).ColumnWidth  = $wsSeq.Columns.Item(9).ColumnWidth
$wsInt.Columns.Item(10).ColumnWidth = $wsSeq.Columns.Item(10).ColumnWidth

# ---------------------------------------------------------------------
# 3. Rename the sheets (do this before touching chart formulas so the
#    SERIES() text below already uses the final sheet names).
# ---------------------------------------------------------------------
$wsDouble.Name = "Single double"
$wsSeq.Name    = "Sequence of items"
$wsInt.Name    = "Single int"

# ---------------------------------------------------------------------
# 4. Chart 1 ("Format single double") stays on sheet1: just repoint its
#    series formulas at the renamed sheet.
# ---------------------------------------------------------------------
$chart1 = $wsDouble.ChartObjects().Item(1).Chart
$chart1.SeriesCollection(1).Formula = "=SERIES('Single double'!`$D`$4,'Single double'!`$E`$3:`$J`$3,'Single double'!`$E`$4:`$J`$4,1)"
$chart1.SeriesCollection(2).Formula = "=SERIES('Single double'!`$D`$5,'Single double'!`$E`$3:`$J`$3,'Single double'!`$E`$5:`$J`$5,2)"
$chart1.SeriesCollection(3).Formula = "=SERIES('Single double'!`$D`$6,'Single double'!`$E`$3:`$J`$3,'Single double'!`$E`$6:`$J`$6,3)"
$chart1.SeriesCollection(4).Formula = "=SERIES('Single double'!`$D`$7,'Single double'!`$E`$3:`$J`$3,'Single double'!`$E`$7:`$J`$7,4)"

# ---------------------------------------------------------------------
# 5. Chart 2 ("Format sequence of several items") used to live on sheet1
#    pointing at rows 36-40; it now belongs on the "Sequence of items"
#    sheet, pointing at rows 3-7. The host has no working "move chart to
#    another sheet" primitive, so re-create it there and drop the old one.
# ---------------------------------------------------------------------
$oldSeqChart = $wsDouble.ChartObjects().Item(2)
$oldSeqChart.Delete()

$seqCo = $wsSeq.ChartObjects().Add(277, 182, 514, 328)
$seqChart = $seqCo.Chart
$seqChart.ChartType = 51
$seqChart.HasTitle = $true
$seqChart.ChartTitle.Text = "Format sequence of several items`n(1000000 iterations)"
$seqChart.HasLegend = $true
$seqChart.Legend.Position = -4107

$s1 = $seqChart.SeriesCollection().NewSeries()
$s1.Formula = "=SERIES('Sequence of items'!`$D`$4,'Sequence of items'!`$E`$3:`$J`$3,'Sequence of items'!`$E`$4:`$J`$4,1)"
$s2 = $seqChart.SeriesCollection().NewSeries()
$s2.Formula = "=SERIES('Sequence of items'!`$D`$5,'Sequence of items'!`$E`$3:`$J`$3,'Sequence of items'!`$E`$5:`$J`$5,2)"
$s3 = $seqChart.SeriesCollection().NewSeries()
$s3.Formula = "=SERIES('Sequence of items'!`$D`$6,'Sequence of items'!`$E`$3:`$J`$3,'Sequence of items'!`$E`$6:`$J`$6,3)"
$s4 = $seqChart.SeriesCollection().NewSeries()
$s4.Formula = "=SERIES('Sequence of items'!`$D`$7,'Sequence of items'!`$E`$3:`$J`$3,'Sequence of items'!`$E`$7:`$J`$7,4)"

$seqValAx = $seqChart.Axes(2)
$seqValAx.HasTitle = $true
$seqValAx.AxisTitle.Text = "Measured time [s]"
$seqValAx.HasMajorGridlines = $true
$seqChart.Axes(1).HasMajorGridlines = $true

# ---------------------------------------------------------------------
# 6. Brand-new Chart 3 ("Format single int") on the "Single int" sheet.
# ---------------------------------------------------------------------
$intCo = $wsInt.ChartObjects().Add(277, 182, 514, 308)
$intChart = $intCo.Chart
$intChart.ChartType = 51
$intChart.HasTitle = $true
$intChart.ChartTitle.Text = "Format single int`n(10000000 iterations)"
$intChart.HasLegend = $true
$intChart.Legend.Position = -4107

$i1 = $intChart.SeriesCollection().NewSeries()
$i1.Formula = "=SERIES('Single int'!`$D`$4,'Single int'!`$E`$3:`$J`$3,'Single int'!`$E`$4:`$J`$4,1)"
$i2 = $intChart.SeriesCollection().NewSeries()
$i2.Formula = "=SERIES('Single int'!`$D`$5,'Single int'!`$E`$3:`$J`$3,'Single int'!`$E`$5:`$J`$5,2)"
$i3 = $intChart.SeriesCollection().NewSeries()
$i3.Formula = "=SERIES('Single int'!`$D`$6,'Single int'!`$E`$3:`$J`$3,'Single int'!`$E`$6:`$J`$6,3)"
$i4 = $intChart.SeriesCollection().NewSeries()
$i4.Formula = "=SERIES('Single int'!`$D`$7,'Single int'!`$E`$3:`$J`$3,'Single int'!`$E`$7:`$J`$7,4)"

$intValAx = $intChart.Axes(2)
$intValAx.HasTitle = $true
$intValAx.AxisTitle.Text = "Measured time [s]"
$intValAx.HasMajorGridlines = $true
$intChart.Axes(1).HasMajorGridlines = $true

# ---------------------------------------------------------------------
# 7. Selections / active sheet, set last so the final Activate() wins.
# ---------------------------------------------------------------------
$wsDouble.Range("K3").Select()
$wsSeq.Columns.Item(11).Select()
$wsInt.Range("P7:P10").Select()
$wsInt.Activate()
